# Add season record columns (Wins, Losses, Ties) to the sheet.
# The original scrape only pulled team/player statistics; this adds the
# season record (Wins/Losses/Ties) as three new columns (AD, AE, AF) that
# apply to every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the number of data rows from the used range (currently A1:AC43)
$lastRow = $ws.UsedRange.Rows.Count

# Header row (row 1) - set the labels and copy the header style used by
# the rest of row 1 (bold, bordered, centered).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the season record for every data row.
$wins = 84
$losses = 78
$ties = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = $wins    # column AD = 30
    $ws.Cells.Item($r, 31).Value = $losses  # column AE = 31
    $ws.Cells.Item($r, 32).Value = $ties    # column AF = 32
}
